$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 23, shifting existing rows 23..160 down to 24..161
$ws.Rows("23:23").Insert()

# Populate the newly inserted row 23 with the new record's data.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R mirror the constant values used throughout
# this sheet (same market/category); only D,J,K,L,M,P differ for this entry.
$ws.Range("A23").Value = 4
$ws.Range("B23").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C23").Value = "Los Lagos"
$ws.Range("D23").Value = 44462
$ws.Range("E23").Value = 10
$ws.Range("F23").Value = 100112037
$ws.Range("G23").Value = "Cebollín"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 90
$ws.Range("K23").Value = 6500
$ws.Range("L23").Value = 6500
$ws.Range("M23").Value = 6500
$ws.Range("N23").Value = "$/paquete 36 unidades"
$ws.Range("O23").Value = "Región Metropolitana"
$ws.Range("P23").Value = 181
$ws.Range("Q23").Value = 36
$ws.Range("R23").Value = "Hortaliza"
